$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-09-09 10:14:45"
$wsZhCn.Range("H4").Value = "2016-09-09 10:14:33"
$wsZhCn.Range("K4").Value = "2016-09-09 10:15:37"
$wsDeDe.Range("H4").Value = "2016-09-09 10:14:45"
$wsDeDe.Range("K4").Value = "2016-09-09 10:15:55"
